$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row => Date (with dashes), D, E, F, G, H
$data = @(
    @{ Row = 3;  Date = "28-07-2022"; D = 2; E = 0; F = 0; G = 2; H = 1 },
    @{ Row = 4;  Date = "01-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 5;  Date = "04-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 6;  Date = "08-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 7;  Date = "11-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 8;  Date = "15-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 9;  Date = "18-08-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 10; Date = "22-08-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 11; Date = "25-08-2022"; D = 2; E = 1; F = 1; G = 0; H = 0 },
    @{ Row = 12; Date = "29-08-2022"; D = 2; E = 1; F = 1; G = 0; H = 0 },
    @{ Row = 13; Date = "01-09-2022"; D = 2; E = 1; F = 1; G = 0; H = 0 },
    @{ Row = 14; Date = "05-09-2022"; D = 2; E = 1; F = 1; G = 0; H = 0 },
    @{ Row = 15; Date = "08-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 16; Date = "12-09-2022"; D = 1; E = 1; F = 0; G = 0; H = 0 },
    @{ Row = 17; Date = "15-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 18; Date = "19-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 19; Date = "22-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 20; Date = "26-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 },
    @{ Row = 21; Date = "29-09-2022"; D = 0; E = 0; F = 0; G = 0; H = 1 }
)

foreach ($entry in $data) {
    $r = $entry.Row
    $cellA = $ws.Range("A$r")
    $cellA.NumberFormat = "@"
    $cellA.Value = $entry.Date
    $cellA.ClearFormats()
    $ws.Range("D$r").Value = $entry.D
    $ws.Range("E$r").Value = $entry.E
    $ws.Range("F$r").Value = $entry.F
    $ws.Range("G$r").Value = $entry.G
    $ws.Range("H$r").Value = $entry.H
}
